# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the order of "Santa Lucia" and "Timor Oriental" in the country list.
# Row 202 currently shows "Santa Lucia", row 203 currently shows "Timor Oriental".
# After the edit they should be swapped (Timor Oriental before Santa Lucia).
$a202 = $ws.Range("A202").Value()
$a203 = $ws.Range("A203").Value()
$ws.Range("A202").Value = $a203
$ws.Range("A203").Value = $a202

# --- Update the "last updated" timestamp text (A1).
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 15:16"

# --- Update statistic values for the countries that changed.

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5033723
$ws.Range("C4").Value = 1544
$ws.Range("D4").Value = 2577914
$ws.Range("E4").Value = 2292973
$ws.Range("G4").Value = 32
$ws.Range("H4").Value = 162836

# India (row 6)
$ws.Range("B6").Value = 2035337
$ws.Range("C6").Value = 9928
$ws.Range("D6").Value = 1382471
$ws.Range("E6").Value = 611153
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = 41713

# Arabia Saudita (row 16)
$ws.Range("B16").Value = 285793
$ws.Range("C16").Value = 1567
$ws.Range("D16").Value = 248948
$ws.Range("E16").Value = 33752
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = 3093

# Irak (row 24)
$ws.Range("B24").Value = 144064
$ws.Range("C24").Value = 3461
$ws.Range("D24").Value = 103197
$ws.Range("E24").Value = 35631
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = 5236

# Suecia (row 34)
$ws.Range("B34").Value = 82323
$ws.Range("C34").Value = 41
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 5763

# Kuwait (row 41)
$ws.Range("B41").Value = 70727
$ws.Range("C41").Value = 682
$ws.Range("D41").Value = 62330
$ws.Range("E41").Value = 7926
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 471

# Paises Bajos (row 45)
$ws.Range("B45").Value = 57501
$ws.Range("C45").Value = 519
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 6154

# Afganistan (row 57)
$ws.Range("B57").Value = 37015
$ws.Range("C57").Value = 119
$ws.Range("E57").Value = 9805
$ws.Range("G57").Value = 9
$ws.Range("H57").Value = 1307

# Azerbaiyan (row 60)
$ws.Range("B60").Value = 33376
$ws.Range("C60").Value = 129
$ws.Range("D60").Value = 29696
$ws.Range("E60").Value = 3197
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 483

# Dinamarca (row 78)
$ws.Range("B78").Value = 14442
$ws.Range("C78").Value = 136
$ws.Range("D78").Value = 12840
$ws.Range("E78").Value = 985

# Croacia (row 101)
$ws.Range("B101").Value = 5466
$ws.Range("C101").Value = 62
$ws.Range("D101").Value = 4758
$ws.Range("E101").Value = 553
